# Auto-generated edit script applying scheduled runner updates to Sheets/Belias_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 116
$ws.Cells.Item(116, 8).Value = 7220.7  # H116: 7998.375 -> 7220.7
$ws.Cells.Item(116, 9).Value = 2457.1428  # I116: 1462 -> 2457.1428
$ws.Cells.Item(116, 10).Value = 18335.666  # J116: 11920.2 -> 18335.666
$ws.Cells.Item(116, 11).Value = 2457.1428  # K116: 1462 -> 2457.1428
$ws.Cells.Item(116, 12).Value = 18335.666  # L116: 11920.2 -> 18335.666
$ws.Cells.Item(116, 13).Value = 984.8571999999999  # M116: 1980 -> 984.8571999999999
$ws.Cells.Item(116, 14).Value = -25219.666  # N116: -18804.2 -> -25219.666

# Row 132
$ws.Cells.Item(132, 8).Value = 15874706  # H132: 17095820 -> 15874706
$ws.Cells.Item(132, 9).Value = 1950133.8  # I132: 2137624.5 -> 1950133.8
$ws.Cells.Item(132, 11).Value = 5850401.4  # K132: 6412873.5 -> 5850401.4
$ws.Cells.Item(132, 13).Value = -5847871.4  # M132: -6410343.5 -> -5847871.4

# Row 137
$ws.Cells.Item(137, 8).Value = 1872.5588  # H137: 1897.0938 -> 1872.5588
$ws.Cells.Item(137, 9).Value = 1274.0869  # I137: 1286.091 -> 1274.0869
$ws.Cells.Item(137, 10).Value = 3123.9092  # J137: 3241.3 -> 3123.9092
$ws.Cells.Item(137, 11).Value = 3822.2607  # K137: 3858.273 -> 3822.2607
$ws.Cells.Item(137, 12).Value = 9371.7276  # L137: 9723.900000000001 -> 9371.7276
$ws.Cells.Item(137, 13).Value = -1272.2607  # M137: -1308.273 -> -1272.2607
$ws.Cells.Item(137, 14).Value = -14471.7276  # N137: -14823.9 -> -14471.7276

$ws = $wb.Worksheets.Item("ARM")
# Row 23
$ws.Cells.Item(23, 8).Value = 20816.5  # H23: 19271.285 -> 20816.5
$ws.Cells.Item(23, 10).Value = 20816.5  # J23: 19271.285 -> 20816.5
$ws.Cells.Item(23, 12).Value = 20816.5  # L23: 19271.285 -> 20816.5
$ws.Cells.Item(23, 14).Value = -21334.5  # N23: -19789.285 -> -21334.5

# Row 32
$ws.Cells.Item(32, 8).Value = 4279.073  # H32: 4169.3623 -> 4279.073
$ws.Cells.Item(32, 9).Value = 2925.5774  # I32: 2791.1572 -> 2925.5774
$ws.Cells.Item(32, 10).Value = 13015.272  # J32: 13816.8 -> 13015.272
$ws.Cells.Item(32, 11).Value = 2925.5774  # K32: 2791.1572 -> 2925.5774
$ws.Cells.Item(32, 12).Value = 13015.272  # L32: 13816.8 -> 13015.272
$ws.Cells.Item(32, 13).Value = -2638.5774  # M32: -2504.1572 -> -2638.5774
$ws.Cells.Item(32, 14).Value = -13589.272  # N32: -14390.8 -> -13589.272

# Row 137
$ws.Cells.Item(137, 8).Value = 41666.668  # H137: 43000 -> 41666.668
$ws.Cells.Item(137, 10).Value = 42200  # J137: 43800 -> 42200
$ws.Cells.Item(137, 12).Value = 42200  # L137: 43800 -> 42200
$ws.Cells.Item(137, 14).Value = -52400  # N137: -54000 -> -52400

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 957.1  # H16: 954.55554 -> 957.1
$ws.Cells.Item(16, 10).Value = 1041.5714  # J16: 1051.8334 -> 1041.5714
$ws.Cells.Item(16, 12).Value = 1041.5714  # L16: 1051.8334 -> 1041.5714
$ws.Cells.Item(16, 14).Value = -1615.5714  # N16: -1625.8334 -> -1615.5714

# Row 31
$ws.Cells.Item(31, 8).Value = 20001616  # H31: 20001532 -> 20001616
$ws.Cells.Item(31, 9).Value = 31250950  # I31: 27778700 -> 31250950
$ws.Cells.Item(31, 10).Value = 2803  # J31: 3103.1428 -> 2803
$ws.Cells.Item(31, 11).Value = 31250950  # K31: 27778700 -> 31250950
$ws.Cells.Item(31, 12).Value = 2803  # L31: 3103.1428 -> 2803
$ws.Cells.Item(31, 13).Value = -31250655  # M31: -27778405 -> -31250655
$ws.Cells.Item(31, 14).Value = -3393  # N31: -3693.1428 -> -3393

# Row 34
$ws.Cells.Item(34, 8).Value = 20001616  # H34: 20001532 -> 20001616
$ws.Cells.Item(34, 9).Value = 31250950  # I34: 27778700 -> 31250950
$ws.Cells.Item(34, 10).Value = 2803  # J34: 3103.1428 -> 2803
$ws.Cells.Item(34, 11).Value = 31250950  # K34: 27778700 -> 31250950
$ws.Cells.Item(34, 12).Value = 2803  # L34: 3103.1428 -> 2803
$ws.Cells.Item(34, 13).Value = -31250748  # M34: -27778498 -> -31250748
$ws.Cells.Item(34, 14).Value = -3207  # N34: -3507.1428 -> -3207

# Row 109
$ws.Cells.Item(109, 8).Value = 26998  # H109: 25998 -> 26998
$ws.Cells.Item(109, 10).Value = 26998  # J109: 25998 -> 26998
$ws.Cells.Item(109, 12).Value = 26998  # L109: 25998 -> 26998
$ws.Cells.Item(109, 14).Value = -29078  # N109: -28078 -> -29078

# Row 113
$ws.Cells.Item(113, 8).Value = 957.1  # H113: 954.55554 -> 957.1
$ws.Cells.Item(113, 10).Value = 1041.5714  # J113: 1051.8334 -> 1041.5714
$ws.Cells.Item(113, 12).Value = 1041.5714  # L113: 1051.8334 -> 1041.5714
$ws.Cells.Item(113, 14).Value = -5381.5714  # N113: -5391.8334 -> -5381.5714

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Cells.Item(34, 8).Value = 388.75  # H34: 362.2 -> 388.75
$ws.Cells.Item(34, 10).Value = 433  # J34: 0 -> 433
$ws.Cells.Item(34, 12).Value = 1299  # L34: 0 -> 1299
$ws.Cells.Item(34, 14).Value = -1467  # N34: None -> -1467

# Row 39
$ws.Cells.Item(39, 8).Value = 3216.5  # H39: 1877.2727 -> 3216.5
$ws.Cells.Item(39, 9).Value = 900  # I39: 550 -> 900
$ws.Cells.Item(39, 10).Value = 3679.8  # J39: 2983.3333 -> 3679.8
$ws.Cells.Item(39, 11).Value = 2700  # K39: 1650 -> 2700
$ws.Cells.Item(39, 12).Value = 11039.4  # L39: 8949.999899999999 -> 11039.4
$ws.Cells.Item(39, 13).Value = -2406  # M39: -1356 -> -2406
$ws.Cells.Item(39, 14).Value = -11627.4  # N39: -9537.999899999999 -> -11627.4

# Row 55
$ws.Cells.Item(55, 8).Value = 597  # H55: 610 -> 597
$ws.Cells.Item(55, 10).Value = 621.5  # J55: 650 -> 621.5
$ws.Cells.Item(55, 12).Value = 1864.5  # L55: 1950 -> 1864.5
$ws.Cells.Item(55, 14).Value = -2218.5  # N55: -2304 -> -2218.5

# Row 75
$ws.Cells.Item(75, 8).Value = 2459.0557  # H75: 4223.091 -> 2459.0557
$ws.Cells.Item(75, 9).Value = 1566.3334  # I75: 2000 -> 1566.3334
$ws.Cells.Item(75, 10).Value = 2637.6  # J75: 4717.1113 -> 2637.6
$ws.Cells.Item(75, 11).Value = 4699.0002  # K75: 6000 -> 4699.0002
$ws.Cells.Item(75, 12).Value = 7912.799999999999  # L75: 14151.3339 -> 7912.799999999999
$ws.Cells.Item(75, 13).Value = -3701.0002  # M75: -5002 -> -3701.0002
$ws.Cells.Item(75, 14).Value = -9908.799999999999  # N75: -16147.3339 -> -9908.799999999999

# Row 78
$ws.Cells.Item(78, 8).Value = 2459.0557  # H78: 4223.091 -> 2459.0557
$ws.Cells.Item(78, 9).Value = 1566.3334  # I78: 2000 -> 1566.3334
$ws.Cells.Item(78, 10).Value = 2637.6  # J78: 4717.1113 -> 2637.6
$ws.Cells.Item(78, 11).Value = 14097.0006  # K78: 18000 -> 14097.0006
$ws.Cells.Item(78, 12).Value = 23738.4  # L78: 42454.00169999999 -> 23738.4
$ws.Cells.Item(78, 13).Value = -9105.000599999999  # M78: -13008 -> -9105.000599999999
$ws.Cells.Item(78, 14).Value = -33722.39999999999  # N78: -52438.00169999999 -> -33722.39999999999

# Row 100
$ws.Cells.Item(100, 8).Value = 3854.5454  # H100: 3950 -> 3854.5454
$ws.Cells.Item(100, 10).Value = 3854.5454  # J100: 3950 -> 3854.5454
$ws.Cells.Item(100, 12).Value = 11563.6362  # L100: 11850 -> 11563.6362
$ws.Cells.Item(100, 14).Value = -13185.6362  # N100: -13472 -> -13185.6362

# Row 110
$ws.Cells.Item(110, 8).Value = 4027  # H110: 1740.875 -> 4027
$ws.Cells.Item(110, 9).Value = 4027  # I110: 1740.875 -> 4027
$ws.Cells.Item(110, 11).Value = 12081  # K110: 5222.625 -> 12081
$ws.Cells.Item(110, 13).Value = -7991  # M110: -1132.625 -> -7991

# Row 113
$ws.Cells.Item(113, 8).Value = 1894471.8  # H113: 1783061.2 -> 1894471.8
$ws.Cells.Item(113, 9).Value = 3030821.8  # I113: 2525766 -> 3030821.8
$ws.Cells.Item(113, 10).Value = 555.3333  # J113: 569.8 -> 555.3333
$ws.Cells.Item(113, 11).Value = 9092465.399999999  # K113: 7577298 -> 9092465.399999999
$ws.Cells.Item(113, 12).Value = 1665.9999  # L113: 1709.4 -> 1665.9999
$ws.Cells.Item(113, 13).Value = -9090295.399999999  # M113: -7575128 -> -9090295.399999999
$ws.Cells.Item(113, 14).Value = -6005.9999  # N113: -6049.4 -> -6005.9999

# Row 115
$ws.Cells.Item(115, 8).Value = 1680.2727  # H115: 2055.9285 -> 1680.2727
$ws.Cells.Item(115, 10).Value = 2925  # J115: 3142.8572 -> 2925
$ws.Cells.Item(115, 12).Value = 8775  # L115: 9428.571599999999 -> 8775
$ws.Cells.Item(115, 14).Value = -11125  # N115: -11778.5716 -> -11125

# Row 120
$ws.Cells.Item(120, 8).Value = 5055.4287  # H120: 5346 -> 5055.4287
$ws.Cells.Item(120, 9).Value = 5055.4287  # I120: 5346 -> 5055.4287
$ws.Cells.Item(120, 11).Value = 15166.2861  # K120: 16038 -> 15166.2861
$ws.Cells.Item(120, 13).Value = -10328.2861  # M120: -11200 -> -10328.2861

# Row 131
$ws.Cells.Item(131, 8).Value = 863.25  # H131: 870.62 -> 863.25
$ws.Cells.Item(131, 9).Value = 532.4167  # I131: 553.5454999999999 -> 532.4167
$ws.Cells.Item(131, 10).Value = 908.36365  # J131: 909.80896 -> 908.36365
$ws.Cells.Item(131, 11).Value = 1597.2501  # K131: 1660.6365 -> 1597.2501
$ws.Cells.Item(131, 12).Value = 2725.09095  # L131: 2729.42688 -> 2725.09095
$ws.Cells.Item(131, 13).Value = 3442.7499  # M131: 3379.3635 -> 3442.7499
$ws.Cells.Item(131, 14).Value = -12805.09095  # N131: -12809.42688 -> -12805.09095

# Row 132
$ws.Cells.Item(132, 8).Value = 3952916.5  # H132: 4235181.5 -> 3952916.5
$ws.Cells.Item(132, 9).Value = 1803906.5  # I132: 1906918.2 -> 1803906.5
$ws.Cells.Item(132, 10).Value = 13892087  # J132: 15876499 -> 13892087
$ws.Cells.Item(132, 11).Value = 16235158.5  # K132: 17162263.8 -> 16235158.5
$ws.Cells.Item(132, 12).Value = 125028783  # L132: 142888491 -> 125028783
$ws.Cells.Item(132, 13).Value = -16232628.5  # M132: -17159733.8 -> -16232628.5
$ws.Cells.Item(132, 14).Value = -125033843  # N132: -142893551 -> -125033843

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Cells.Item(57, 8).Value = 14315.23  # H57: 15972.546 -> 14315.23
$ws.Cells.Item(57, 10).Value = 15383.167  # J57: 17419.8 -> 15383.167
$ws.Cells.Item(57, 12).Value = 15383.167  # L57: 17419.8 -> 15383.167
$ws.Cells.Item(57, 14).Value = -17023.167  # N57: -19059.8 -> -17023.167

# Row 110
$ws.Cells.Item(110, 8).Value = 33400  # H110: 32800 -> 33400
$ws.Cells.Item(110, 10).Value = 33400  # J110: 32800 -> 33400
$ws.Cells.Item(110, 12).Value = 33400  # L110: 32800 -> 33400
$ws.Cells.Item(110, 14).Value = -41580  # N110: -40980 -> -41580

# Row 114
$ws.Cells.Item(114, 8).Value = 32000  # H114: 28699.5 -> 32000
$ws.Cells.Item(114, 10).Value = 32000  # J114: 28699.5 -> 32000
$ws.Cells.Item(114, 12).Value = 32000  # L114: 28699.5 -> 32000
$ws.Cells.Item(114, 14).Value = -40678  # N114: -37377.5 -> -40678

# Row 116
$ws.Cells.Item(116, 8).Value = 30000  # H116: 0 -> 30000
$ws.Cells.Item(116, 10).Value = 30000  # J116: 0 -> 30000
$ws.Cells.Item(116, 12).Value = 30000  # L116: 0 -> 30000
$ws.Cells.Item(116, 14).Value = -39178  # N116: None -> -39178

# Row 119
$ws.Cells.Item(119, 8).Value = 32000  # H119: 31142.857 -> 32000
$ws.Cells.Item(119, 10).Value = 32000  # J119: 31142.857 -> 32000
$ws.Cells.Item(119, 12).Value = 32000  # L119: 31142.857 -> 32000
$ws.Cells.Item(119, 14).Value = -41676  # N119: -40818.857 -> -41676

# Row 122
$ws.Cells.Item(122, 8).Value = 204190  # H122: 103635 -> 204190
$ws.Cells.Item(122, 9).Value = 252737.5  # I122: 114038.89 -> 252737.5
$ws.Cells.Item(122, 11).Value = 758212.5  # K122: 342116.67 -> 758212.5
$ws.Cells.Item(122, 13).Value = -755762.5  # M122: -339666.67 -> -755762.5

# Row 124
$ws.Cells.Item(124, 8).Value = 26400  # H124: 25250 -> 26400
$ws.Cells.Item(124, 10).Value = 26400  # J124: 25250 -> 26400
$ws.Cells.Item(124, 12).Value = 26400  # L124: 25250 -> 26400
$ws.Cells.Item(124, 14).Value = -36220  # N124: -35070 -> -36220

# Row 126
$ws.Cells.Item(126, 8).Value = 1717.25  # H126: 2157 -> 1717.25
$ws.Cells.Item(126, 9).Value = 1554  # I126: 2000 -> 1554
$ws.Cells.Item(126, 10).Value = 2207  # J126: 2209.3333 -> 2207
$ws.Cells.Item(126, 11).Value = 4662  # K126: 6000 -> 4662
$ws.Cells.Item(126, 12).Value = 6621  # L126: 6627.999899999999 -> 6621
$ws.Cells.Item(126, 13).Value = -2192  # M126: -3530 -> -2192
$ws.Cells.Item(126, 14).Value = -11561  # N126: -11567.9999 -> -11561

# Row 128
$ws.Cells.Item(128, 8).Value = 0  # H128: 20000 -> 0
$ws.Cells.Item(128, 10).Value = 0  # J128: 20000 -> 0
$ws.Cells.Item(128, 12).Value = 0  # L128: 20000 -> 0
$ws.Cells.Item(128, 14).ClearContents()  # N128: delete (was -29960)

# Row 129
$ws.Cells.Item(129, 8).Value = 49999  # H129: 25000 -> 49999
$ws.Cells.Item(129, 9).Value = 0  # I129: 20000 -> 0
$ws.Cells.Item(129, 10).Value = 49999  # J129: 30000 -> 49999
$ws.Cells.Item(129, 11).Value = 0  # K129: 20000 -> 0
$ws.Cells.Item(129, 12).Value = 49999  # L129: 30000 -> 49999
$ws.Cells.Item(129, 13).ClearContents()  # M129: delete (was -15000)
$ws.Cells.Item(129, 14).Value = -59999  # N129: -40000 -> -59999

# Row 130
$ws.Cells.Item(130, 8).Value = 1000000  # H130: 118000 -> 1000000
$ws.Cells.Item(130, 10).Value = 1000000  # J130: 118000 -> 1000000
$ws.Cells.Item(130, 12).Value = 1000000  # L130: 118000 -> 1000000
$ws.Cells.Item(130, 14).Value = -1010040  # N130: -128040 -> -1010040

$ws = $wb.Worksheets.Item("LTW")
# Row 105
$ws.Cells.Item(105, 8).Value = 29000  # H105: 29500 -> 29000
$ws.Cells.Item(105, 10).Value = 29000  # J105: 29500 -> 29000
$ws.Cells.Item(105, 12).Value = 29000  # L105: 29500 -> 29000
$ws.Cells.Item(105, 14).Value = -35988  # N105: -36488 -> -35988

# Row 108
$ws.Cells.Item(108, 8).Value = 29713  # H108: 29836 -> 29713
$ws.Cells.Item(108, 10).Value = 29713  # J108: 29836 -> 29713
$ws.Cells.Item(108, 12).Value = 29713  # L108: 29836 -> 29713
$ws.Cells.Item(108, 14).Value = -37393  # N108: -37516 -> -37393

# Row 110
$ws.Cells.Item(110, 8).Value = 22224.5  # H110: 21849.625 -> 22224.5
$ws.Cells.Item(110, 10).Value = 22224.5  # J110: 21849.625 -> 22224.5
$ws.Cells.Item(110, 12).Value = 22224.5  # L110: 21849.625 -> 22224.5
$ws.Cells.Item(110, 14).Value = -30404.5  # N110: -30029.625 -> -30404.5

# Row 114
$ws.Cells.Item(114, 8).Value = 0  # H114: 29900 -> 0
$ws.Cells.Item(114, 10).Value = 0  # J114: 29900 -> 0
$ws.Cells.Item(114, 12).Value = 0  # L114: 29900 -> 0
$ws.Cells.Item(114, 14).ClearContents()  # N114: delete (was -38578)

# Row 132
$ws.Cells.Item(132, 8).Value = 2658.0635  # H132: 2837.1475 -> 2658.0635
$ws.Cells.Item(132, 9).Value = 2525.0205  # I132: 2751.7874 -> 2525.0205
$ws.Cells.Item(132, 11).Value = 7575.0615  # K132: 8255.3622 -> 7575.0615
$ws.Cells.Item(132, 13).Value = -5045.0615  # M132: -5725.3622 -> -5045.0615

$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Cells.Item(61, 8).Value = 15111.4  # H61: 15213.75 -> 15111.4
$ws.Cells.Item(61, 9).Value = 7375  # I61: 8991.833000000001 -> 7375
$ws.Cells.Item(61, 10).Value = 20269  # J61: 21435.666 -> 20269
$ws.Cells.Item(61, 11).Value = 7375  # K61: 8991.833000000001 -> 7375
$ws.Cells.Item(61, 12).Value = 20269  # L61: 21435.666 -> 20269
$ws.Cells.Item(61, 13).Value = -7083  # M61: -8699.833000000001 -> -7083
$ws.Cells.Item(61, 14).Value = -20853  # N61: -22019.666 -> -20853

# Row 110
$ws.Cells.Item(110, 8).Value = 0  # H110: 29800 -> 0
$ws.Cells.Item(110, 10).Value = 0  # J110: 29800 -> 0
$ws.Cells.Item(110, 12).Value = 0  # L110: 29800 -> 0
$ws.Cells.Item(110, 14).ClearContents()  # N110: delete (was -37980)

# Row 114
$ws.Cells.Item(114, 8).Value = 25000  # H114: 28333.334 -> 25000
$ws.Cells.Item(114, 10).Value = 25000  # J114: 28333.334 -> 25000
$ws.Cells.Item(114, 12).Value = 25000  # L114: 28333.334 -> 25000
$ws.Cells.Item(114, 14).Value = -33678  # N114: -37011.334 -> -33678

# Row 116
$ws.Cells.Item(116, 8).Value = 31600  # H116: 31000 -> 31600
$ws.Cells.Item(116, 10).Value = 31600  # J116: 31000 -> 31600
$ws.Cells.Item(116, 12).Value = 31600  # L116: 31000 -> 31600
$ws.Cells.Item(116, 14).Value = -40778  # N116: -40178 -> -40778

# Row 119
$ws.Cells.Item(119, 8).Value = 32344.908  # H119: 30828.215 -> 32344.908
$ws.Cells.Item(119, 10).Value = 32344.908  # J119: 30828.215 -> 32344.908
$ws.Cells.Item(119, 12).Value = 32344.908  # L119: 30828.215 -> 32344.908
$ws.Cells.Item(119, 14).Value = -42020.908  # N119: -40504.215 -> -42020.908

# Row 123
$ws.Cells.Item(123, 8).Value = 25000  # H123: 21658.334 -> 25000
$ws.Cells.Item(123, 10).Value = 25000  # J123: 21658.334 -> 25000
$ws.Cells.Item(123, 12).Value = 25000  # L123: 21658.334 -> 25000
$ws.Cells.Item(123, 14).Value = -34800  # N123: -31458.334 -> -34800
